# Apply "backup functional hicpro 2.11.1" update to extra_data workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A (type) for the two new rows
$ws.Range("A102").Value = "hic"
$ws.Range("A103").Value = "longrange"

# Column B (name) for the two new rows
$ws.Range("B102").Value = "JIA_merged_HiC"
$ws.Range("B103").Value = "JIA_merged_pCHiC"

# Column C (url) for the two new rows
$ws.Range("C102").Value = "http://bartzabel.ls.manchester.ac.uk/worthingtonlab/psa_functional_genomics/JIA_CHiC/jia_hic_merged_MBOI.allValidPairs.hic"
$ws.Range("C103").Value = "http://bartzabel.ls.manchester.ac.uk/worthingtonlab/psa_functional_genomics/JIA_CHiC/jia_merged_washU_text.txt.new_washu.bed.gz"

# Column D (sample) for the two new rows
$ws.Range("D102").Value = "JIA_CHiC"
$ws.Range("D103").Value = "JIA_CHiC"

# Match the styling used by the block of rows above (A: name-style font, D: date-format style)
$ws.Range("A99").Copy()
$ws.Range("A102:A103").PasteSpecial(-4122)

$ws.Range("D99").Copy()
$ws.Range("D102:D103").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update selection/view similar to the recorded commit
$ws.Range("C106").Select()
$excel.ActiveWindow.ScrollRow = 73
